# Adds the three new assignment worksheets (10, 11, 12) required for V3,
# matching the header/row layout of the existing "Assignment N" test sheets.
#
# Notes on this interpreter's quirks (discovered empirically):
#  - COM worksheet/range objects must NOT be passed as function parameters;
#    doing so silently loses the live reference. Helper functions instead
#    take a sheet index and re-resolve $wb.Worksheets.Item($idx) internally.
#  - Named parameters (-Foo bar) are not bound correctly; use positional
#    parameters everywhere instead.
#  - A parenthesized expression used as a non-first positional argument to
#    a function call (e.g. `Foo 1 ($c + 1)`) gets mis-parsed. Always
#    pre-compute such expressions into a temp variable first.

$wb = $excel.ActiveWorkbook

function Set-CellValue($idx, $r, $c, $v) {
    $ws = $wb.Worksheets.Item($idx)
    $ws.Cells.Item($r, $c).Value = $v
}

function Style-HeaderRow($idx, $lastCol) {
    $ws = $wb.Worksheets.Item($idx)
    $headerRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item(1, $lastCol))
    $headerRange.Font.Bold = $true
    $headerRange.Borders.LineStyle = 1
    $headerRange.HorizontalAlignment = -4108
    $headerRange.VerticalAlignment = -4160
}

function Add-AssignmentSheet($SheetName, $Headers, $Rows) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newWs = $wb.Worksheets.Add($null, $lastSheet)
    $newWs.Name = $SheetName
    $idx = $wb.Worksheets.Count

    $headerCount = $Headers.Length
    for ($c = 0; $c -lt $headerCount; $c++) {
        $col = $c + 1
        Set-CellValue $idx 1 $col $Headers[$c]
    }
    Style-HeaderRow $idx $headerCount

    $rowCount = $Rows.Length
    for ($r = 0; $r -lt $rowCount; $r++) {
        $row = $Rows[$r]
        $sheetRow = $r + 2
        $colCount = $row.Length
        for ($c = 0; $c -lt $colCount; $c++) {
            $val = $row[$c]
            $col = $c + 1
            if ($null -ne $val) {
                Set-CellValue $idx $sheetRow $col $val
            }
        }
    }
}

# ---------------------------------------------------------------------------
# Assignment 10 - Func Test
# ---------------------------------------------------------------------------
$headers10 = @("test_type", "function_name", "description", "solution_file", "test_inputs", "tolerance")
$rows10 = @(
    ,@("function_exists", "calculate_stats", "Function calculate_stats should exist", $null, $null, $null)
    ,@("test_function_solution", "calculate_stats", "Test function with lists and arrays against solution", "solutions/assignment10_solution.py", "[{'args': [[1, 2, 3, 4, 5]]}, {'args': [[10, 20, 30]]}, {'args': [np.array([5, 10, 15])]}]", 0.01)
    ,@("function_not_called", "np.mean", "Should NOT use np.mean - must calculate manually", $null, $null, $null)
    ,@("function_not_called", "numpy.mean", "Should NOT use numpy.mean", $null, $null, $null)
)
Add-AssignmentSheet "Assignment 10 - Func Test" $headers10 $rows10

# ---------------------------------------------------------------------------
# Assignment 11 - Relations
# ---------------------------------------------------------------------------
$headers11 = @("test_type", "variable_name", "expected_value", "tolerance", "description", "var1_name", "var2_name", "relationship")
$rows11 = @(
    ,@("variable_value", "x", "[0, 1, 2, 3, 4, 5]", 0, "x should be array of values", $null, $null, $null)
    ,@("check_relationship", $null, $null, 0.001, "y should equal cos(π * x)", "x", "y", "lambda x: np.cos(np.pi * x)")
    ,@("check_relationship", $null, $null, 0.001, "z should equal 2x + 1", "x", "z", "lambda x: 2*x + 1")
    ,@("variable_type", "y", "list", $null, "y should be a list or array", $null, $null, $null)
)
Add-AssignmentSheet "Assignment 11 - Relations" $headers11 $rows11

# ---------------------------------------------------------------------------
# Assignment 12 - Adv Plot
# ---------------------------------------------------------------------------
$headers12 = @("test_type", "description", "min_lines", "function", "min_length", "tolerance", "title", "xlabel", "ylabel", "has_legend", "has_grid", "function_name")
$rows12 = @(
    ,@("plot_created", "Should create a plot", $null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    ,@("check_multiple_lines", "Plot should have at least 2 lines", 2, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    ,@("check_function_any_line", "One line should be cos(2x) with at least 50 points", $null, "lambda x: np.cos(2*x)", 50, 0.01, $null, $null, $null, $null, $null, $null)
    ,@("plot_properties", "Plot should have proper labels and legend", $null, $null, $null, $null, "Trigonometric Functions", "x", "y", "'true", "'true", $null)
    ,@("function_not_called", "Should NOT use np.linspace", $null, $null, $null, $null, $null, $null, $null, $null, $null, "np.linspace")
    ,@("function_not_called", "Should NOT use numpy.linspace", $null, $null, $null, $null, $null, $null, $null, $null, $null, "numpy.linspace")
)
Add-AssignmentSheet "Assignment 12 - Adv Plot" $headers12 $rows12
